$d = $word.ActiveDocument

function Get-ParagraphByExactText($doc, $exactText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        $t = $p.Range.Text
        if ($t -eq $exactText -or $t -eq ($exactText + "`r")) {
            return $p
        }
    }
    throw "Paragraph not found for text: $exactText"
}

function Get-RunXml($runText) {
    # Mirror Word's own serialisation: only stamp xml:space="preserve" when
    # the text has leading/trailing whitespace that would otherwise be lost.
    if ($runText -ne $runText.Trim()) {
        return '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">' + $runText + '</w:t></w:r>'
    }
    return '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>' + $runText + '</w:t></w:r>'
}

function Set-ParagraphTwoRuns($doc, $paragraph, $run1Text, $run2Text) {
    # Range covering the paragraph's text, excluding the trailing paragraph mark.
    $full = $paragraph.Range
    $start = $full.Start
    $end = $full.End
    $textRange = $doc.Range($start, $end - 1)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
        (Get-RunXml $run1Text) + (Get-RunXml $run2Text) + `
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $textRange.InsertXML($xml)
}

function Set-ParagraphOneRun($doc, $paragraph, $runText) {
    $full = $paragraph.Range
    $start = $full.Start
    $end = $full.End
    $textRange = $doc.Range($start, $end - 1)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
        (Get-RunXml $runText) + `
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $textRange.InsertXML($xml)
}

# 1) Numbered list of evaluation criteria: "- " -> "1) " / "2)" / "3)" / "4)"
$p1 = Get-ParagraphByExactText $d "- L’outil doit être facile à utiliser."
Set-ParagraphTwoRuns $d $p1 "1) " "L’outil doit être facile à utiliser."

$p2 = Get-ParagraphByExactText $d "- L’outil doit être accessible pour tout le monde."
Set-ParagraphTwoRuns $d $p2 "2)" " L’outil doit être accessible pour tout le monde."

$p3 = Get-ParagraphByExactText $d "- Le développeur doit facilement modifier son projet dans l’outil avec le moins d’étapes possible."
Set-ParagraphTwoRuns $d $p3 "3)" " Le développeur doit facilement modifier son projet dans l’outil avec le moins d’étapes possible."

$p4 = Get-ParagraphByExactText $d "- Le coût de l’utilisation de l’outil doit être le plus bas possible, préférablement gratuit."
Set-ParagraphTwoRuns $d $p4 "4)" " Le coût de l’utilisation de l’outil doit être le plus bas possible, préférablement gratuit."

# 2) Merge the three runs ("... retourner de" + "s" + " «commits» en arrière.") into one run.
$p5 = Get-ParagraphByExactText $d "Tu as aussi la possibilité de créer plusieurs branches de développement, ce qui te permet de répartir plusieurs étapes de développement sans avoir à retourner des «commits» en arrière."
Set-ParagraphOneRun $d $p5 "Tu as aussi la possibilité de créer plusieurs branches de développement, ce qui te permet de répartir plusieurs étapes de développement sans avoir à retourner des «commits» en arrière."
